# Replace the M2Doc field (a Word field whose code is the query
# "m:'doc.html'.fromHTMLURI()") with the same text typed out literally,
# wrapped in curly braces, e.g. "{m:'doc.html'.fromHTMLURI()}".
#
# This mirrors what TokenIteratorFieldRewriterSplit expects: plain text
# runs instead of a real field (fldChar begin/instrText.../fldChar end),
# while keeping the _GoBack bookmark exactly where it was, sitting
# between the "doc.html" run and the "'.fromHTMLURI()" run.

$d = $word.ActiveDocument

# Locate the paragraph that holds the M2Doc field.
$fieldPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $fieldPara = $p
    }
}

$r = $fieldPara.Range

# Build a WordProcessingML package fragment that replaces the whole
# paragraph with literal text runs (one run per former instrText run,
# minus the leading/trailing " " instrText runs which are absorbed into
# the new "{" and "}" runs), preserving the bookmark in place.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' +
       '<w:r><w:t>{</w:t></w:r>' +
       '<w:r><w:t>m</w:t></w:r>' +
       '<w:r><w:t>:</w:t></w:r>' +
       '<w:r><w:t>''</w:t></w:r>' +
       '<w:r><w:t>doc.html</w:t></w:r>' +
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
       '<w:bookmarkEnd w:id="0"/>' +
       '<w:r><w:t>''.fromHTMLURI()</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
